# Encode categorical MDF channels numerically with NV enums.
#
# Adds a "Type" column (D) to the vbRcSignals sheet that records each
# signal's value type: "float" for the plain numeric channels and
# "enum" for the channels whose categorical/state values are now
# encoded numerically via NV enums.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vbRcSignals")

# Header cell - copy C1's formatting (bold header style) onto D1.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "Type"

# Data rows 2-16 share the plain data-row formatting used by columns A:C.
$ws.Range("C2").Copy()
$ws.Range("D2:D16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Numeric (float) signals: rows 2-14.
$ws.Range("D2:D14").Value = "float"

# Categorical / state signals now represented with an NV enum: rows 15-20
# (aebFullState, aebPartialState, isVehStoppedNV, fcwRequest, aebRequest,
# dbsRequest). Rows 17-20 keep the workbook's default (unstyled) look,
# same as columns A/C already do on those rows.
$ws.Range("D15:D20").Value = "enum"

# Restore the selected cell to match the saved view.
$ws.Range("F14").Select()
